# Refresh the cryptocurrency price table from the latest coinranking.com
# scrape (scheduled GitHub Actions update).
#
# Column D ("Price") holds numeric-looking values in a cell that carries no
# explicit number format, so a plain numeric-looking string would silently
# be reinterpreted by Excel as a real number (dropping e.g. the trailing
# zero in "0.990"). Typing a leading apostrophe is what forces Excel to
# store the literal text instead - exactly as a user would do at the
# keyboard - and re-applying the "Normal" cell style afterwards clears the
# transient quote-prefix flag picked up along the way so the cell's
# formatting stays byte-for-byte what it was before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText "D2" "68.132.76"
$ws.Range("E2").Value = "  -6.49%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.699.20"
$ws.Range("E3").Value = "  -6.63%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "0.990"
$ws.Range("E4").Value = "  -1.04%  "

# Row 5 - BNB
Set-PriceText "D5" "579.83"
$ws.Range("E5").Value = "  -5.07%  "

# Row 6 - Solana
Set-PriceText "D6" "175.51"
$ws.Range("E6").Value = "  +5.27%  "

# Row 7 - LidoStakedEther
Set-PriceText "D7" "3.828.05"
$ws.Range("E7").Value = "  -3.22%  "

# Row 8 - XRP
Set-PriceText "D8" "0.634"
$ws.Range("E8").Value = "  -6.43%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.54%  "

# Row 10 - Cardano
Set-PriceText "D10" "0.718"
$ws.Range("E10").Value = "  -4.22%  "

# Row 11 - Dogecoin
Set-PriceText "D11" "0.167"
$ws.Range("E11").Value = "  -8.32%  "

# Row 12 - Avalanche
Set-PriceText "D12" "52.57"
$ws.Range("E12").Value = "  -6.36%  "

# Row 13 - ShibaInu
Set-PriceText "D13" "0.0000304"
$ws.Range("E13").Value = "  -8.99%  "

# Row 14 - Polkadot
Set-PriceText "D14" "10.71"
$ws.Range("E14").Value = "  -2.95%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-PriceText "D15" "4.252.75"
$ws.Range("E15").Value = "  -7.66%  "

# Row 16 - WrappedEther
Set-PriceText "D16" "3.686.23"
$ws.Range("E16").Value = "  -7.25%  "

# Row 17 - was TRON, becomes Chainlink
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-PriceText "D17" "19.51"
$ws.Range("E17").Value = "  -4.35%  "

# Row 18 - was Chainlink, becomes TRON
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-PriceText "D18" "0.127"
$ws.Range("E18").Value = "  -2.93%  "

# Row 19 - Uniswap
Set-PriceText "D19" "13.05"
$ws.Range("E19").Value = "  -6.83%  "

# Row 20 - Polygon
Set-PriceText "D20" "1.14"
$ws.Range("E20").Value = "  -7.69%  "

# Row 21 - WrappedBTC
Set-PriceText "D21" "67.345.07"
$ws.Range("E21").Value = "  -7.49%  "

# Row 22 - BitcoinCash
Set-PriceText "D22" "411.36"
$ws.Range("E22").Value = "  -5.81%  "

# Row 23 - PancakeSwap
Set-PriceText "D23" "4.57"
$ws.Range("E23").Value = "  -6.10%  "

# Row 24 - Litecoin
Set-PriceText "D24" "88.56"
$ws.Range("E24").Value = "  -7.03%  "

# Row 25 - ImmutableX
Set-PriceText "D25" "3.09"
$ws.Range("E25").Value = "  -7.94%  "

# Row 26 - InternetComputer(DFINITY)
Set-PriceText "D26" "12.88"
$ws.Range("E26").Value = "  -8.79%  "

# Row 27 - RenderToken
Set-PriceText "D27" "10.82"
$ws.Range("E27").Value = "  -2.07%  "

# Row 28 - Toncoin
Set-PriceText "D28" "3.82"
$ws.Range("E28").Value = "  -6.52%  "

# Row 29 - LEO
Set-PriceText "D29" "5.94"
$ws.Range("E29").Value = "  -0.29%  "

# Row 30 - Filecoin
Set-PriceText "D30" "9.64"
$ws.Range("E30").Value = "  -7.24%  "

# Row 31 - NEARProtocol
Set-PriceText "D31" "8.20"
$ws.Range("E31").Value = "  +3.79%  "

# Row 32 - EthereumClassic
Set-PriceText "D32" "33.07"
$ws.Range("E32").Value = "  -7.83%  "

# Row 33 - Cosmos
Set-PriceText "D33" "12.82"
$ws.Range("E33").Value = "  -5.82%  "

# Row 34 - OKB
Set-PriceText "D34" "68.29"
$ws.Range("E34").Value = "  -3.26%  "

# Row 35 - was Hedera, becomes InjectiveProtocol
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-PriceText "D35" "44.92"
$ws.Range("E35").Value = "  -5.77%  "

# Row 36 - was InjectiveProtocol, becomes Hedera
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-PriceText "D36" "0.119"
$ws.Range("E36").Value = "  -8.36%  "

# Row 37 - PEPE
Set-PriceText "D37" "0.0₃0937"
$ws.Range("E37").Value = "  -7.20%  "

# Row 38 - Bittensor
Set-PriceText "D38" "591.20"
$ws.Range("E38").Value = "  -7.52%  "

# Row 39 - TheGraph
Set-PriceText "D39" "0.405"
$ws.Range("E39").Value = "  -5.60%  "

# Row 40 - Dai
Set-PriceText "D40" "0.998"
$ws.Range("E40").Value = "  -0.23%  "

# Row 41 - was FirstDigitalUSD, becomes dogwifhat
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-PriceText "D41" "3.23"
$ws.Range("E41").Value = "  +13.42%  "

# Row 42 - was dogwifhat, becomes FirstDigitalUSD
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-PriceText "D42" "0.980"
$ws.Range("E42").Value = "  -2.19%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  -7.73%  "

# Row 44 - ThetaToken
Set-PriceText "D44" "3.09"
$ws.Range("E44").Value = "  -10.66%  "

# Row 45 - VeChain
Set-PriceText "D45" "0.0446"
$ws.Range("E45").Value = "  -7.44%  "

# Row 46 - THORChain
Set-PriceText "D46" "9.59"
$ws.Range("E46").Value = "  -11.30%  "

# Row 47 - Fetch.AI
Set-PriceText "D47" "2.59"
$ws.Range("E47").Value = "  -0.16%  "

# Row 48 - Stellar
Set-PriceText "D48" "0.136"
$ws.Range("E48").Value = "  -8.35%  "

# Row 49 - Maker
Set-PriceText "D49" "2.746.32"
$ws.Range("E49").Value = "  -3.68%  "

# Row 50 - WEMIXToken
Set-PriceText "D50" "2.65"
$ws.Range("E50").Value = "  -17.21%  "

# Row 51 - ApeXProtocol
Set-PriceText "D51" "3.12"
$ws.Range("E51").Value = "  -8.49%  "
